# Apply updated "views/likes" counts (column F) across sheets, as generated
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 13355
$ws1.Range("F7").Value  = 8
$ws1.Range("F9").Value  = 124
$ws1.Range("F11").Value = 67
$ws1.Range("F14").Value = 13333
$ws1.Range("F15").Value = 330
$ws1.Range("F17").Value = 8891
$ws1.Range("F19").Value = 7966
$ws1.Range("F20").Value = 241
$ws1.Range("F22").Value = 138
$ws1.Range("F23").Value = 426
$ws1.Range("F26").Value = 17
$ws1.Range("F27").Value = 1012
$ws1.Range("F29").Value = 18
$ws1.Range("F32").Value = 152
$ws1.Range("F33").Value = 368

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 31

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 13355
$ws4.Range("F8").Value  = 8
$ws4.Range("F10").Value = 124
$ws4.Range("F12").Value = 67
$ws4.Range("F15").Value = 13334
$ws4.Range("F16").Value = 330
$ws4.Range("F18").Value = 8891
$ws4.Range("F20").Value = 7966
$ws4.Range("F21").Value = 241
$ws4.Range("F23").Value = 138
$ws4.Range("F24").Value = 426
$ws4.Range("F27").Value = 17
$ws4.Range("F28").Value = 1012
$ws4.Range("F30").Value = 18
$ws4.Range("F31").Value = 31
$ws4.Range("F35").Value = 152
$ws4.Range("F36").Value = 368
